# ENH: Optional initialization through steady-state estimation
#
# Adds a new "steady" column (O) to the "vars" sheet. This flags, per
# variable, whether its initial value should be estimated by running the
# model to steady state instead of being taken straight from "default".
# Antibiotic (A) starts at a fixed (non-steady) concentration, while the
# biological state variables (D, R, T) are initialized via steady-state
# estimation.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("vars")

# New header cell, formatted like the other header cells (A1:K1).
$ws.Range("O1").Value = "steady"
$ws.Range("A1").Copy()
$ws.Range("O1").PasteSpecial(-4122)

# New data values: A -> 0 (not steady-state init), D/R/T -> 1 (steady-state init).
$ws.Range("O2").Value = 0
$ws.Range("O3").Value = 1
$ws.Range("O4").Value = 1
$ws.Range("O5").Value = 1

# Match the plain data-cell formatting used by the rest of the table.
$ws.Range("D2").Copy()
$ws.Range("O2:O5").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Keep the sheet's active selection in sync with the new right-most column.
$ws.Range("O6").Select()

# Cosmetic: the tab-area/scrollbar split ratio also shrank in the source
# commit; carried over here for parity (best effort - window chrome, not data).
$excel.ActiveWindow.TabRatio = 0.255
